$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.05%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.65%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.810"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.25%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06326"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.16%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.920"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.27%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.364"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.21%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8830"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.78%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9464"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.16%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1478"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.23%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05266"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.03%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07346"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.61%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03151"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.18%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09069"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001549"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.28%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006278"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.77%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005798"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.08%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.465"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.41%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.277"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "6.08%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.76%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1342"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.77%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.906"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-4.84%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04322"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.90%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001177"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.18%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003586"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-11.67%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001693"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-12.62%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04017"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.47%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006644"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "58.64%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1165"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.71%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002341"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "16.41%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.00%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005240"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.56%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "822.16%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02254"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "6.26%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.11%"
